$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix BOM designators for R1 and R3:
#  - R1 belongs with R6 (row 13)
#  - R3 belongs with R2, R7, R8 (row 15)
# (leading apostrophe keeps these as quote-prefixed text, matching the
#  original cell formatting for the Designator column)
$ws.Range("C13").Value = "'R1,R6"
$ws.Range("C15").Value = "'R2, R3, R7, R8"

# Remove the now-redundant placeholder row (old row 16: NC/NC/"R1, R3"/0805R/NC/NC)
$ws.Rows("16").Delete()

# Widen column E (Part) so the merged designator text is readable
$ws.Columns("E").ColumnWidth = 14.666666666666666

# Update selection to match the saved state
$ws.Range("A16:XFD16").Select()
